$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'257.06"
$ws.Range("E2").Value = "'4.76%"
$ws.Range("D3").Value = "'27.46"
$ws.Range("E3").Value = "'-3.20%"
$ws.Range("D4").Value = "'5.206"
$ws.Range("E4").Value = "'-1.61%"
$ws.Range("D5").Value = "'0.05910"
$ws.Range("E5").Value = "'3.47%"
$ws.Range("E6").Value = "'0.38%"
$ws.Range("D7").Value = "'0.8639"
$ws.Range("E7").Value = "'1.62%"
$ws.Range("D8").Value = "'1.003"
$ws.Range("E8").Value = "'13.51%"
$ws.Range("B9").Value = "'One"
$ws.Range("C9").Value = "'https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D9").Value = "'0.01048"
$ws.Range("E9").Value = "'1,652.88%"
$ws.Range("B10").Value = "'WazirX"
$ws.Range("C10").Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1418"
$ws.Range("E10").Value = "'1.40%"
$ws.Range("D11").Value = "'0.07185"
$ws.Range("E11").Value = "'1.41%"
$ws.Range("D12").Value = "'0.03150"
$ws.Range("E12").Value = "'0.13%"
$ws.Range("D13").Value = "'0.09216"
$ws.Range("E13").Value = "'-0.11%"
$ws.Range("E14").Value = "'0.43%"
$ws.Range("B15").Value = "'TigerCash"
$ws.Range("C15").Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005870"
$ws.Range("E15").Value = "'-1.11%"
$ws.Range("B16").Value = "'LEO"
$ws.Range("C16").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.493"
$ws.Range("E16").Value = "'-0.10%"
$ws.Range("B17").Value = "'GateToken"
$ws.Range("C17").Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'3.268"
$ws.Range("E17").Value = "'1.88%"
$ws.Range("B18").Value = "'BTSEToken"
$ws.Range("C18").Value = "'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "'2.224"
$ws.Range("E18").Value = "'1.64%"
$ws.Range("B19").Value = "'BitpandaEcosystemToken"
$ws.Range("C19").Value = "'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "'0.3172"
$ws.Range("E19").Value = "'0.17%"
$ws.Range("B20").Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range("C20").Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D20").Value = "'0.03575"
$ws.Range("E20").Value = "'7.40%"
$ws.Range("D21").Value = "'0.1306"
$ws.Range("E21").Value = "'0.95%"
$ws.Range("D22").Value = "'3.542"
$ws.Range("E22").Value = "'0.82%"
$ws.Range("D23").Value = "'0.04185"
$ws.Range("E23").Value = "'2.78%"
$ws.Range("E24").Value = "'1.49%"
$ws.Range("D25").Value = "'0.001217"
$ws.Range("E25").Value = "'-0.57%"
$ws.Range("D26").Value = "'0.004514"
$ws.Range("E26").Value = "'8.66%"
$ws.Range("E27").Value = "'0.04%"
$ws.Range("E28").Value = "'34.15%"
$ws.Range("D40").Value = "'0.03817"
$ws.Range("E40").Value = "'0.91%"
$ws.Range("D41").Value = "'0.005610"
$ws.Range("E41").Value = "'49.78%"
$ws.Range("E42").Value = "'3.16%"
$ws.Range("D43").Value = "'0.001800"
$ws.Range("E43").Value = "'-20.67%"
$ws.Range("E44").Value = "'12.45%"
$ws.Range("D45").Value = "'0.00005438"
$ws.Range("E45").Value = "'3.08%"
$ws.Range("E46").Value = "'0.03%"
$ws.Range("D48").Value = "'0.002178"
$ws.Range("E48").Value = "'-4.01%"
$ws.Range("E49").Value = "'0.03%"
$ws.Range("E50").Value = "'0.03%"
